# Populate the "Lương" (Salary) worksheet with the salary-calculation
# report rows: column A = the line-item label, column B = its amount.
# All amounts start at 0 — this sheet is the new calculation base that
# downstream code fills in, replacing the old (manually generated) report.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")
$ws.Activate()

$rows = @(
    ,@("Danh mục", 0)
    ,@("Ngày công", 0)
    ,@("Phụ cấp", 0)
    ,@("Lương cơ bản tại CẦN THƠ", 0)
    ,@("Chiết khấu sale chính tại CẦN THƠ", 0)
    ,@("Chiết khấu sale phụ tại CẦN THƠ", 0)
    ,@("Đơn 1 bác sĩ tại CẦN THƠ", 0)
    ,@("Đơn 2 bác sĩ tại CẦN THƠ", 0)
    ,@("Công phụ phẫu 1 tại CẦN THƠ", 0)
    ,@("Công phụ phẫu 2 tại CẦN THƠ", 0)
    ,@("Lương cơ bản tại LONG XUYÊN", 0)
    ,@("Chiết khấu sale chính tại LONG XUYÊN", 0)
    ,@("Chiết khấu sale phụ tại LONG XUYÊN", 0)
    ,@("Đơn 1 bác sĩ tại LONG XUYÊN", 0)
    ,@("Đơn 2 bác sĩ tại LONG XUYÊN", 0)
    ,@("Công phụ phẫu 1 tại LONG XUYÊN", 0)
    ,@("Công phụ phẫu 2 tại LONG XUYÊN", 0)
    ,@("Lương cơ bản tại SÓC TRĂNG", 0)
    ,@("Chiết khấu sale chính tại SÓC TRĂNG", 0)
    ,@("Chiết khấu sale phụ tại SÓC TRĂNG", 0)
    ,@("Đơn 1 bác sĩ tại SÓC TRĂNG", 0)
    ,@("Đơn 2 bác sĩ tại SÓC TRĂNG", 0)
    ,@("Công phụ phẫu 1 tại SÓC TRĂNG", 0)
    ,@("Công phụ phẫu 2 tại SÓC TRĂNG", 0)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
}
